$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dim")
$ws.Name = "sql πεδία"
